$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

function Set-CellText($cell, $text) {
    $cell.Range.Text = $text
    $cell.Range.Font.Name = "Calibri"
}

$row6 = $t.Rows.Item(6)
Set-CellText $row6.Cells.Item(1) "19/07/2024"
Set-CellText $row6.Cells.Item(2) "report"
Set-CellText $row6.Cells.Item(3) "Chapter 1: introduction completed and uploaded"

$row7 = $t.Rows.Item(7)
Set-CellText $row7.Cells.Item(1) "31/07/2024"
Set-CellText $row7.Cells.Item(2) "report"
Set-CellText $row7.Cells.Item(3) "Chapter 2: literature review completed and uploaded"
